$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for two new rows. The existing row 2 ("3204 1800 00" /
#    "Carotenoid colouring matters...") shifts down to become row 4; the
#    newly inserted rows 2 and 3 will hold the two new entries.
# ---------------------------------------------------------------------------
$ws.Range("A2:A3").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Header row text tweak (singular -> plural).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "10 Digit Codes Introduced from 1 January 2022"
$ws.Range("B1").Value = "Code Description"

# ---------------------------------------------------------------------------
# 3. New row 2: masterbatch pigment mixture code (this becomes the umbrella
#    heading's detail row further down visually, but in sheet order it is
#    row 2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "                                3204 1800 00"
$ws.Range("B2").Value = "Carotenoid colouring matters and preparations based thereon"

# ---------------------------------------------------------------------------
# 4. New row 3: umbrella subheading / "Other" row.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "                 Umbrella Subheading 2903 4970"
$ws.Range("B3").Value = "Other"

# ---------------------------------------------------------------------------
# 5. Row 4 (former row 2): new code + long multi-line description.
# ---------------------------------------------------------------------------
$desc = "Concentrated mixture of pigments (masterbatch) in the form of pellets containing by weight:" + [char]10 + `
  "-50 % or more but not more than 70 % of polyamide-6.6 (CAS RN 32131-17-2)," + [char]10 + `
  "-15 % or more but not more than 20 % of iron powder (CAS RN 7439-89-6)," + [char]10 + `
  "-5 % or more but not more than 15 % of barium sulphate (CAS RN 7727-43-7), and" + [char]10 + `
  "-5 % or more but not more than 10 % of blue pigment, consisting of a mixture of Titanium dioxide (CAS RN 13463-67-7) and Copper(II) phtalocyanine(CAS RN 147-14-8)"

$ws.Range("A4").Value = "                                3206 4970 50"
$ws.Range("B4").Value = $desc

# ---------------------------------------------------------------------------
# 6. Formatting.
#
#    A2 / B2 / A3 / B3 all share one plain style: Arial 9, default colour
#    (theme text colour), vertically centred. The existing B4 cell (the old
#    row-2 "Carotenoid..." cell, style index 3 pre-edit) already carries that
#    font, so copy its format across rather than touching Font.* directly
#    (direct Font property sets clone a brand-new font entry every time and
#    would bloat the style table far beyond what the source file has).
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy() | Out-Null
$ws.Range("A2:B3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# A4 already carries the Arial-9 / no-colour font that used to live on the
# pre-edit A2 cell, so there's no font work to do here - just strip the
# wrap/number-format baggage that cell carried and add the thin right-hand
# border + left/centre alignment the new layout calls for.
$ws.Range("A4").NumberFormat = "general"
$ws.Range("A4").WrapText = $false
$ws.Range("A4").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A4").VerticalAlignment = -4108     # xlCenter
$ws.Range("A4").Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous

# B4 keeps its inherited font (Arial 9, theme colour, vertical-centred) but
# becomes a wrapped "text" cell to hold the long multi-line description.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4108     # xlCenter

# Row 4 needs extra height to show the five-line description.
$ws.Rows.Item(4).RowHeight = 69.75

# ---------------------------------------------------------------------------
# 7. Selection matches the post-edit cursor position recorded in the diff.
# ---------------------------------------------------------------------------
$ws.Range("B3").Select() | Out-Null
